# Updated data sheet and code of further pages
# Adds new lead-data columns (T:AC) to Sheet1's header/value rows, applies
# the matching number formats, widens the new columns, and refreshes the
# active-cell selections on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: new header cells (row 1) + new value cells (row 2) -----------
# Entered column-group by column-group (T,U,V together, then W, then X, ...)
# so the shared-string table is populated in the same interleaved order as
# the source data entry.
$ws1.Range("T1").Value = "ifsc"
$ws1.Range("U1").Value = "Branch"
$ws1.Range("V1").Value = "relType"

$ws1.Range("T2").Value = "hdfc0001111"
$ws1.Range("U2").Value = "hdfc"
$ws1.Range("V2").Value = "Self"

$ws1.Range("W1").Value = "acc"

$ws1.Range("X1").Value = "fathername"
$ws1.Range("X2").Value = "abc"

$ws1.Range("Y1").Value = "motherName"
$ws1.Range("Y2").Value = "def"

$ws1.Range("Z1").Value = "wieght"

$ws1.Range("AA2").Value = " 3I INFOTECH "
$ws1.Range("AA1").Value = "org"

$ws1.Range("AB1").Value = "natureOfWork"
$ws1.Range("AB2").Value = "Web Developers"

$ws1.Range("AC1").Value = "dobb"

# Numeric cells with their custom formats
$ws1.Range("W2").Value = 89674523456789
$ws1.Range("W2").NumberFormat = "# ?/?"
$ws1.Range("W2").HorizontalAlignment = -4131

$ws1.Range("Z2").Value = 56

$ws1.Range("AC2").Value = 36212
$ws1.Range("AC2").NumberFormat = "mm-dd-yy"

# --- Column widths for the newly used columns ------------------------------
$ws1.Columns.Item(17).ColumnWidth = 14.72       # Q
$ws1.Columns.Item(20).ColumnWidth = 11.72       # T
$ws1.Columns.Item(23).ColumnWidth = 20.39       # W
$ws1.Columns.Item(27).ColumnWidth = 16.28       # AA
$ws1.Columns.Item(28).ColumnWidth = 22.95       # AB
$ws1.Columns.Item(29).ColumnWidth = 9.5         # AC

# --- Selections -------------------------------------------------------------
$ws1.Range("AC3").Select()

$ws2.Range("A1:D2").Select()
